# ---------------------------------------------------------------
# edit.ps1 - apply "New crime data collected" update to 61st Precinct
# CompStat weekly report (rows 14-30 data block), plus the header
# volume/week-range shared strings and the column-E autofit width.
# ---------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ------------------------
# "Volume 31   Number  36" -> "...Number  37"
$ws.Range("C1").Characters(32, 2).Text = "37"

# "Report Covering the Week  9/2/2024  Through  9/8/2024"
#   -> "...Week  9/9/2024  Through  9/15/2024"
$ws.Range("C9").Characters(28, 7).Text = "9/9/2024"
$ws.Range("C9").Characters(45, 8).Text = "9/15/2024"

# --- Helper: convert a cell to a NUMBER style by copying number ---
# format from a stable donor cell elsewhere on the sheet, then set
# the numeric value (two-step paste keeps cellXfs/styles.xml stable)
function Set-AsNumber($addr, $donor, $val) {
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($addr).Value = $val
}

# --- Helper: convert a cell to the TEXT "N/A" style by copying ----
# format+value from a stable donor cell that already holds the same
# literal text (e.g. "0" or "***.*") as a shared string.
function Set-AsText($addr, $donor) {
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# --- Row 14 ---
$ws.Range("N14").Value = -66.666666666666

# --- Row 15 ---
Set-AsNumber "C15" "I31" 1
Set-AsNumber "F15" "I31" 1
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 11
$ws.Range("K15").Value = -8.333333333333
$ws.Range("L15").Value = 83.333333333333
$ws.Range("M15").Value = 266.666666666667
$ws.Range("N15").Value = -35.294117647058

# --- Row 16 ---
Set-AsText "C16" "D14"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 104
$ws.Range("J16").Value = 113
$ws.Range("K16").Value = -7.964601769911
$ws.Range("L16").Value = -20.610687022900
$ws.Range("M16").Value = -2.803738317757
$ws.Range("N16").Value = -83.492063492063

# --- Row 17 ---
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -50
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = -41.666666666666
$ws.Range("I17").Value = 187
$ws.Range("J17").Value = 166
$ws.Range("K17").Value = 12.650602409638
$ws.Range("L17").Value = 3.314917127071
$ws.Range("M17").Value = 156.164383561644
$ws.Range("N17").Value = -17.256637168141

# --- Row 18 ---
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -50
$ws.Range("G18").Value = 34
$ws.Range("H18").Value = -67.647058823529
$ws.Range("I18").Value = 114
$ws.Range("J18").Value = 163
$ws.Range("K18").Value = -30.061349693251
$ws.Range("L18").Value = -30.487804878048
$ws.Range("M18").Value = -38.378378378378
$ws.Range("N18").Value = -92.271186440678

# --- Row 19 ---
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -18.181818181818
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = -15.217391304347
$ws.Range("I19").Value = 434
$ws.Range("J19").Value = 476
$ws.Range("K19").Value = -8.823529411764
$ws.Range("L19").Value = -20.073664825046
$ws.Range("M19").Value = 65.019011406844
$ws.Range("N19").Value = -57.986447241045

# --- Row 20 ---
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = 4.761904761904
$ws.Range("I20").Value = 176
$ws.Range("J20").Value = 143
$ws.Range("K20").Value = 23.076923076923
$ws.Range("L20").Value = 28.467153284671
$ws.Range("M20").Value = 27.536231884058
$ws.Range("N20").Value = -92.669720949604

# --- Row 21 ---
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -14.285714285714
$ws.Range("F21").Value = 99
$ws.Range("G21").Value = 143
$ws.Range("H21").Value = -30.769230769230
$ws.Range("I21").Value = 1027
$ws.Range("J21").Value = 1077
$ws.Range("K21").Value = -4.642525533890
$ws.Range("L21").Value = -11.769759450171
$ws.Range("M21").Value = 33.031088082901
$ws.Range("N21").Value = -82.247191011236

# --- Row 22 ---
Set-AsText "C22" "D14"
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("L22").Value = 114.285714285714
$ws.Range("M22").Value = 36.363636363636

# --- Row 23 ---
$ws.Range("C23").Value = 1
Set-AsText "D23" "D14"
Set-AsText "E23" "E14"
$ws.Range("F23").Value = 4
$ws.Range("H23").Value = -60
$ws.Range("I23").Value = 51
$ws.Range("K23").Value = 10.869565217391
$ws.Range("L23").Value = -26.086956521739
$ws.Range("M23").Value = 168.421052631579

# --- Row 24 ---
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = -25
$ws.Range("F24").Value = 107
$ws.Range("H24").Value = -2.727272727272
$ws.Range("I24").Value = 834
$ws.Range("J24").Value = 882
$ws.Range("K24").Value = -5.442176870748
$ws.Range("L24").Value = -26.519823788546
$ws.Range("M24").Value = 10.904255319148

# --- Row 25 ---
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 28.571428571428
$ws.Range("F25").Value = 55
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = 89.655172413793
$ws.Range("I25").Value = 355
$ws.Range("J25").Value = 362
$ws.Range("K25").Value = -1.933701657458
$ws.Range("L25").Value = -43.290734824281

# --- Row 26 ---
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 133.333333333333
$ws.Range("F26").Value = 37
$ws.Range("G26").Value = 27
$ws.Range("H26").Value = 37.037037037037
$ws.Range("I26").Value = 333
$ws.Range("J26").Value = 284
$ws.Range("K26").Value = 17.253521126760
$ws.Range("L26").Value = 4.0625
$ws.Range("M26").Value = -5.665722379603

# --- Row 27 ---
Set-AsNumber "C27" "I31" 1
Set-AsNumber "F27" "I31" 1
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 16
$ws.Range("K27").Value = 6.666666666666

# --- Row 28 ---
Set-AsNumber "D28" "I31" 2
Set-AsNumber "E28" "K31" -100
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -33.333333333333
$ws.Range("J28").Value = 30
$ws.Range("K28").Value = -16.666666666666
$ws.Range("L28").Value = -10.714285714285

# --- Row 29 ---
Set-AsNumber "C29" "I31" 2
Set-AsNumber "F29" "I31" 2
$ws.Range("I29").Value = 6
$ws.Range("K29").Value = 100
$ws.Range("L29").Value = -45.454545454545
$ws.Range("M29").Value = 200
$ws.Range("N29").Value = -77.777777777777

# --- Row 30 ---
Set-AsNumber "C30" "I31" 1
Set-AsNumber "F30" "I31" 1
$ws.Range("I30").Value = 4
$ws.Range("K30").Value = 33.333333333333
$ws.Range("L30").Value = -50
$ws.Range("M30").Value = 300
$ws.Range("N30").Value = -80

# --- Column E width (bestFit autofit side-effect of new text) ---
$ws.Columns("E").ColumnWidth = 7.433768
